$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 124
$ws.Range("I2").Value = 284
$ws.Range("J2").Value = 1229
$ws.Range("L2").Value = 291
$ws.Range("M2").Value = 22
$ws.Range("N2").Value = 215
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 16
$ws.Range("S2").Value = 131
$ws.Range("T2").Value = 205
$ws.Range("U2").Value = 22
$ws.Range("V2").Value = 1942
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1941
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 7
